# Update recalculated TPM-based NATMI ligand-receptor metrics for Robo2-Robo2
# (Sheet1, rows 2-10, columns G-T) to reflect the new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 1.302860333333333
$ws.Range("H2").Value = 3.908581
$ws.Range("I2").Value = 0.9669439908960468
$ws.Range("J2").Value = 0.9669439908960467
$ws.Range("M2").Value = 1.302860333333333
$ws.Range("N2").Value = 3.908581
$ws.Range("O2").Value = 0.9669439908960468
$ws.Range("P2").Value = 0.9669439908960467
$ws.Range("Q2").Value = 1.697445048173444
$ws.Range("R2").Value = 15.277005433561
$ws.Range("S2").Value = 0.9349806815299742
$ws.Range("T2").Value = 0.9349806815299739
$ws.Range("G3").Value = 1.302860333333333
$ws.Range("H3").Value = 3.908581
$ws.Range("I3").Value = 0.9669439908960468
$ws.Range("J3").Value = 0.9669439908960467
$ws.Range("O3").Value = 0.008324674682103805
$ws.Range("P3").Value = 0.008324674682103805
$ws.Range("Q3").Value = 0.01461375007222222
$ws.Range("R3").Value = 0.13152375065
$ws.Range("S3").Value = 0.008049494160024732
$ws.Range("T3").Value = 0.008049494160024732
$ws.Range("G4").Value = 1.302860333333333
$ws.Range("H4").Value = 3.908581
$ws.Range("I4").Value = 0.9669439908960468
$ws.Range("J4").Value = 0.9669439908960467
$ws.Range("M4").Value = 0.033323
$ws.Range("N4").Value = 0.099969
$ws.Range("O4").Value = 0.02473133442184949
$ws.Range("P4").Value = 0.02473133442184949
$ws.Range("Q4").Value = 0.04341521488766666
$ws.Range("R4").Value = 0.390736933989
$ws.Range("S4").Value = 0.02391381520604792
$ws.Range("T4").Value = 0.02391381520604792
$ws.Range("I5").Value = 0.008324674682103805
$ws.Range("J5").Value = 0.008324674682103805
$ws.Range("M5").Value = 1.302860333333333
$ws.Range("N5").Value = 3.908581
$ws.Range("O5").Value = 0.9669439908960468
$ws.Range("P5").Value = 0.9669439908960467
$ws.Range("Q5").Value = 0.01461375007222222
$ws.Range("R5").Value = 0.13152375065
$ws.Range("S5").Value = 0.008049494160024732
$ws.Range("T5").Value = 0.008049494160024732
$ws.Range("I6").Value = 0.008324674682103805
$ws.Range("J6").Value = 0.008324674682103805
$ws.Range("O6").Value = 0.008324674682103805
$ws.Range("P6").Value = 0.008324674682103805
$ws.Range("S6").Value = 0.00006930020856286007
$ws.Range("T6").Value = 0.00006930020856286007
$ws.Range("I7").Value = 0.008324674682103805
$ws.Range("J7").Value = 0.008324674682103805
$ws.Range("M7").Value = 0.033323
$ws.Range("N7").Value = 0.099969
$ws.Range("O7").Value = 0.02473133442184949
$ws.Range("P7").Value = 0.02473133442184949
$ws.Range("Q7").Value = 0.0003737729833333333
$ws.Range("R7").Value = 0.00336395685
$ws.Range("S7").Value = 0.0002058803135162127
$ws.Range("T7").Value = 0.0002058803135162127
$ws.Range("G8").Value = 0.033323
$ws.Range("H8").Value = 0.099969
$ws.Range("I8").Value = 0.02473133442184949
$ws.Range("J8").Value = 0.02473133442184949
$ws.Range("M8").Value = 1.302860333333333
$ws.Range("N8").Value = 3.908581
$ws.Range("O8").Value = 0.9669439908960468
$ws.Range("P8").Value = 0.9669439908960467
$ws.Range("Q8").Value = 0.04341521488766666
$ws.Range("R8").Value = 0.390736933989
$ws.Range("S8").Value = 0.02391381520604792
$ws.Range("T8").Value = 0.02391381520604792
$ws.Range("G9").Value = 0.033323
$ws.Range("H9").Value = 0.099969
$ws.Range("I9").Value = 0.02473133442184949
$ws.Range("J9").Value = 0.02473133442184949
$ws.Range("O9").Value = 0.008324674682103805
$ws.Range("P9").Value = 0.008324674682103805
$ws.Range("Q9").Value = 0.0003737729833333333
$ws.Range("R9").Value = 0.00336395685
$ws.Range("S9").Value = 0.0002058803135162127
$ws.Range("T9").Value = 0.0002058803135162127
$ws.Range("G10").Value = 0.033323
$ws.Range("H10").Value = 0.099969
$ws.Range("I10").Value = 0.02473133442184949
$ws.Range("J10").Value = 0.02473133442184949
$ws.Range("M10").Value = 0.033323
$ws.Range("N10").Value = 0.099969
$ws.Range("O10").Value = 0.02473133442184949
$ws.Range("P10").Value = 0.02473133442184949
$ws.Range("Q10").Value = 0.001110422329
$ws.Range("R10").Value = 0.009993800961000001
$ws.Range("S10").Value = 0.0006116389022853573
$ws.Range("T10").Value = 0.0006116389022853573
